$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# B11 currently holds the text "R40" (shared string). The new value is the
# text "1" -- a leading apostrophe forces Excel to store it as text (not a
# number), matching the t="s" shared-string cell produced by the edit.
$ws.Range("B11").Value = "'1"
